$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug in biotime_purchase_types:
# 1) Rename the "product types" input header to the correct name.
$ws.Range("B2").Value = "in:biotime_product_types"

# 2) Fix the intersection-count threshold used by the two "combo package"
#    rules: it should require MORE THAN ONE matching product type (#> 1),
#    not more than zero (#> 0).
$ws.Range("B12").Value = '$in intersection $(Differentiation Kit, Human Embryonic Progenitor Package, Growth Media, Basal Media) #> 1'
$ws.Range("B13").Value = '$in intersection $(Differentiation Kit, Human Embryonic Progenitor Package, Growth Media, Basal Media) #> 1'

# Leave the selection where the edit finished, matching the saved workbook state.
$ws.Range("B14").Select()
